# Architecture and testbench of ALU unit
# Add a new time-tracking entry (row 9) for 27.11.2019 covering the
# architecture & testbench work on the ALU Unit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing time/number formatting from the row above (row 8)
# onto the new row so the new cells render the same way (hh:mm time
# cells for B/C, the shared-formula duration format for D).
$ws.Range("B9").NumberFormat = $ws.Range("B8").NumberFormat
$ws.Range("C9").NumberFormat = $ws.Range("C8").NumberFormat
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat

# New row of data: 27.11.2019, 15:30 - 18:00, Calculator Control Unit /
# ALU Unit, Architecture & Testbench.
$ws.Range("A9").Value = "27.11.2019"
$ws.Range("B9").Value = 0.64583333333333337
$ws.Range("C9").Value = 0.75
$ws.Range("D9").Formula = "=C9-B9"
$ws.Range("E9").Value = "ALU Unit"
$ws.Range("F9").Value = "Architecture, Testbench"

# After entering data across the row, the cursor lands one column past
# the last entry.
[void]$ws.Range("G9").Select()
